$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'Any person, including a trustee, who directs or manages the affairs of an unincorporated organization or association other than a partnership.'
$ws.Range("E2").Value = 0.7473684210526316

$ws.Range("C5").Value = 'A party named in the process, pleadings, or papers served.'

$ws.Range("C6").Value = 'An individual, corporation, partnership, or other unincorporated organization or association that resides in any place not subject to the jurisdiction of the United States.'
$ws.Range("E6").Value = 0.3809523809523809

$ws.Range("C7").Value = 'Legal documents served on a non-resident investment adviser, general partner, or managing agent.'

$ws.Range("C8").Value = ''
$ws.Range("E8").Value = 1.0

$ws.Range("C9").Value = 'Any individual, corporation, partnership, or other unincorporated organization or association.'

$ws.Range("C10").Value = 'An unincorporated organization or association other than a partnership.'

$ws.Range("C12").Value = 'The records maintained by the Commission, which include a copy of the process, pleadings, or papers served.'

$ws.Range("C13").Value = 'The Secretary of the Commission is responsible for forwarding a copy of the process, pleadings, or papers to each named party by registered or certified mail.'

$ws.Range("C14").Value = 'A section of the document describing the service of process, pleadings, or papers on the Commission.'

$ws.Range("C15").Value = 'The main office and place of business, as defined in § 275.203A-3(c) of this chapter.'

$ws.Range("C16").Value = 'A section of the chapter that defines ''Principal office and place of business.'''

$ws.Range("C17").Value = ''
$ws.Range("E17").Value = 1.0

$ws.Range("C18").Value = 'A person who directs or manages the affairs of an unincorporated organization or association.'

$ws.Range("C19").Value = 'A corporation that is incorporated in or has its principal office and place of business in any place not subject to the jurisdiction of the United States.'

$ws.Range("C20").Value = 'An organization or association that is not incorporated and is managed by a managing agent.'

$ws.Range("C21").Value = 'A method of mailing used by the Secretary to forward a copy of the process, pleadings, or papers to each named party.'

$ws.Range("C23").Value = 'The Secretary''s certification that the Commission was served with process, pleadings, or other papers and forwarded these documents to a named party.'

$ws.Range("C24").Value = 'A section of the document describing the forwarding of documents by the Secretary to a named party by registered or certified mail.'

$ws.Range("C25").Value = 'The location where the main business activities are conducted, as defined in § 275.203A-3(c) of this chapter.'

$ws.Range("C26").Value = ''
$ws.Range("E26").Value = 1.0

$ws.Range("C27").Value = 'The activities directed or managed by a managing agent of an unincorporated organization or association.'

$ws.Range("C28").Value = 'Agents appointed to receive service of process, pleadings, or other papers on behalf of a non-resident investment adviser, general partner, or managing agent.'

$ws.Range("C30").Value = 'The legal documents served on a non-resident investment adviser, general partner, or managing agent.'

$ws.Range("C31").Value = ''
$ws.Range("E31").Value = 1.0

$ws.Range("C32").Value = 'The Secretary of the Commission responsible for forwarding documents to named parties.'

$ws.Range("C33").Value = ''
$ws.Range("E33").Value = 1.0

$ws.Range("C34").Value = 'The last address filed with the Commission for a named party.'

$ws.Range("C35").Value = 'A method of mailing used by the Secretary to forward a copy of the process, pleadings, or papers to each named party.'

$ws.Range("C36").Value = 'The entity responsible for receiving service of process, pleadings, or other papers on behalf of non-resident investment advisers, general partners, or managing agents.'

$ws.Range("C38").Value = 'Interested persons must state their reasons for requesting a hearing.'
$ws.Range("E38").Value = 0.1888544891640866

$ws.Range("C39").Value = 'Any section of the Act or any rule or regulation thereunder.'
$ws.Range("E39").Value = 0.8194444444444444

$ws.Range("C40").Value = 'An application for registration as an investment adviser is excluded from the definition of ''application''.'
$ws.Range("E40").Value = 0.5377358490566038

$ws.Range("C42").Value = 'Any person who may submit facts and request a hearing on the matter.'
$ws.Range("E42").Value = 0.2897196261682243

$ws.Range("C44").Value = 'Notice of the initiation of the proceeding will be published in the Federal Register.'
$ws.Range("E44").Value = 1.0

$ws.Range("C47").Value = 'May be ordered if necessary or appropriate in the public interest or for the protection of investors.'
$ws.Range("E47").Value = 0.572289156626506

$ws.Range("C48").Value = 'A hearing may be ordered if necessary or appropriate in the public interest.'
$ws.Range("E48").Value = 0.5114503816793894

$ws.Range("C50").Value = 'Interested persons may submit facts bearing upon the desirability of a hearing.'
$ws.Range("E50").Value = 0.4342857142857143

$ws.Range("C52").Value = 'An ''application'' means any application for an order of the Commission under the Act other than an application for registration as an investment adviser.'

$ws.Range("C54").Value = 'The period within which interested persons may submit facts and request a hearing.'
$ws.Range("E54").Value = 0.3257142857142857

$ws.Range("C55").Value = 'The body that may initiate proceedings and order hearings.'
$ws.Range("E55").Value = 0.2469879518072289

$ws.Range("C57").Value = 'The subject of the proceeding or hearing.'
$ws.Range("E57").Value = 0.2026143790849673

$ws.Range("C58").Value = 'A person is presumed to control a trust if the person is a trustee or managing agent of the trust.'
$ws.Range("E58").Value = 1.0

$ws.Range("C59").Value = 'A small business or small organization for purposes of the Investment Advisers Act of 1940 is an investment adviser with assets under management of less than $25 million, did not have total assets of $5 million or more on the last day of the most recent fiscal year, and does not control, is not controlled by, and is not under common control with another investment adviser with assets under management of $25 million or more.'
$ws.Range("E59").Value = 0.5052410901467506

$ws.Range("C60").Value = 'A person is presumed to control a partnership if the person has the right to receive upon dissolution, or has contributed, 25 percent or more of the capital of the partnership.'
$ws.Range("E60").Value = 1.0

$ws.Range("C61").Value = 'The Investment Advisers Act of 1940 is referenced in defining a small business or small organization as an investment adviser with certain asset criteria.'
$ws.Range("E61").Value = 0.0

$ws.Range("C64").Value = 'A person is presumed to control a trust if the person is a trustee or managing agent of the trust.'
$ws.Range("E64").Value = 1.0

$ws.Range("C65").Value = 'A person is presumed to control a corporation, partnership, LLC, or trust under certain conditions outlined in the document.'
$ws.Range("E65").Value = 0.3241379310344827

$ws.Range("C67").Value = 'A person is presumed to control an LLC if the person has the right to vote 25 percent or more of a class of the interests of the LLC.'
$ws.Range("E67").Value = 0.0

$ws.Range("C68").Value = 'A person is presumed to control a partnership if the person has the right to receive upon dissolution, or has contributed, 25 percent or more of the capital of the partnership.'
$ws.Range("E68").Value = 1.0

$ws.Range("C69").Value = 'Assets under management are defined under Section 203A(a)(3) of the Act and reported on the annual updating amendment to Form ADV.'
$ws.Range("E69").Value = 0.4679245283018868

$ws.Range("C72").Value = 'Total assets means the total assets as shown on the balance sheet of the investment adviser or other person, or the balance sheet of the investment adviser or such other person with its subsidiaries consolidated, whichever is larger.'
$ws.Range("E72").Value = 0.8755364806866953

$ws.Range("C73").Value = 'Total assets means the total assets as shown on the balance sheet of the investment adviser or other person, or the balance sheet of the investment adviser or such other person with its subsidiaries consolidated, whichever is larger.'
$ws.Range("E73").Value = 0.8583690987124464

$ws.Range("C74").Value = 'A small business or small organization for purposes of the Investment Advisers Act of 1940 is an investment adviser with assets under management of less than $25 million, did not have total assets of $5 million or more on the last day of the most recent fiscal year, and does not control, is not controlled by, and is not under common control with another investment adviser with assets under management of $25 million or more.'
$ws.Range("E74").Value = 0.5052410901467506

$ws.Range("C75").Value = 'A person is presumed to control a trust if the person is a trustee or managing agent of the trust.'
$ws.Range("E75").Value = 1.0

$ws.Range("C76").Value = 'A person is presumed to control a corporation if the person directly or indirectly has the right to vote 25 percent or more of a class of the corporation''s voting securities, or has the power to sell or direct the sale of 25 percent or more of a class of the corporation''s voting securities.'
$ws.Range("E76").Value = 0.9965635738831615

$ws.Range("C77").Value = 'An investment adviser is a small business or small organization under the Investment Advisers Act of 1940 if it has assets under management of less than $25 million, did not have total assets of $5 million or more on the last day of the most recent fiscal year, and does not control, is not controlled by, and is not under common control with another investment adviser with assets under management of $25 million or more.'
$ws.Range("E77").Value = 0.5220125786163522

$ws.Range("C78").Value = 'A person is presumed to control a partnership if the person has the right to receive upon dissolution, or has contributed, 25 percent or more of the capital of the partnership.'
$ws.Range("E78").Value = 1.0

$ws.Range("C79").Value = 'A person is presumed to control a corporation if the person has the right to vote 25 percent or more of a class of the corporation''s voting securities.'
$ws.Range("E79").Value = 0.5206896551724138

$ws.Range("C80").Value = 'Control means the power, directly or indirectly, to direct the management or policies of a person, whether through ownership of securities, by contract, or otherwise.'
$ws.Range("E80").Value = 0.9096385542168675

$ws.Range("C82").Value = 'A person is presumed to control a corporation if the person directly or indirectly has the right to vote 25 percent or more of a class of the corporation''s voting securities.'
$ws.Range("E82").Value = 1.0
